# com_protocol_MCU_NavU.xlsx — "Added ICC tester and fixed typo in xlsx"
#
# 1) Fix the CMD-code "typo": on the "Msg NavU -> MCU" sheet the CMD column
#    (B) for the last five commands was off by one (Drive distance repeated
#    Soft stop's 0x03, etc.). Shift each of those five rows to the correct
#    code, which also introduces the new 0x08 code for "Continue Drive".
# 2) Make "Msg NavU -> MCU" the active/selected sheet with C10 selected
#    (previously "Structure" was active/selected).

$wb = $excel.ActiveWorkbook

$wsNavToMcu = $wb.Worksheets.Item("Msg NavU -> MCU")

$wsNavToMcu.Range("B7").Value2  = "0x04"
$wsNavToMcu.Range("B8").Value2  = "0x05"
$wsNavToMcu.Range("B9").Value2  = "0x06"
$wsNavToMcu.Range("B10").Value2 = "0x07"
$wsNavToMcu.Range("B11").Value2 = "0x08"

$wsNavToMcu.Activate()
$wsNavToMcu.Range("C10").Select()
